$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '23.697.14' }
    @{ Cell = 'E2'; Value = '  +1.25%  ' }
    @{ Cell = 'D3'; Value = '1.652.44' }
    @{ Cell = 'E3'; Value = '  +1.11%  ' }
    @{ Cell = 'D4'; Value = '1.003' }
    @{ Cell = 'E4'; Value = '  +0.19%  ' }
    @{ Cell = 'E5'; Value = '  +0.16%  ' }
    @{ Cell = 'D6'; Value = '303.23' }
    @{ Cell = 'E6'; Value = '  -0.05%  ' }
    @{ Cell = 'D7'; Value = '0.3802' }
    @{ Cell = 'E7'; Value = '  +0.55%  ' }
    @{ Cell = 'D8'; Value = '0.3618' }
    @{ Cell = 'E8'; Value = '  -0.04%  ' }
    @{ Cell = 'D9'; Value = '51.08' }
    @{ Cell = 'E9'; Value = '  -1.63%  ' }
    @{ Cell = 'D10'; Value = '1.245' }
    @{ Cell = 'E10'; Value = '  +1.58%  ' }
    @{ Cell = 'D11'; Value = '0.08206' }
    @{ Cell = 'E11'; Value = '  +0.47%  ' }
    @{ Cell = 'D12'; Value = '1.003' }
    @{ Cell = 'E12'; Value = '  +0.16%  ' }
    @{ Cell = 'D13'; Value = '22.62' }
    @{ Cell = 'E13'; Value = '  +1.29%  ' }
    @{ Cell = 'D14'; Value = '6.516' }
    @{ Cell = 'E14'; Value = '  +0.68%  ' }
    @{ Cell = 'D15'; Value = '7.447' }
    @{ Cell = 'E15'; Value = '  +1.29%  ' }
    @{ Cell = 'D16'; Value = '0.00001233' }
    @{ Cell = 'E16'; Value = '  -0.58%  ' }
    @{ Cell = 'D17'; Value = '1.652.92' }
    @{ Cell = 'E17'; Value = '  +1.41%  ' }
    @{ Cell = 'D18'; Value = '97.32' }
    @{ Cell = 'E18'; Value = '  +2.59%  ' }
    @{ Cell = 'D19'; Value = '0.07013' }
    @{ Cell = 'E19'; Value = '  +1.17%  ' }
    @{ Cell = 'D20'; Value = '6.787' }
    @{ Cell = 'E20'; Value = '  +3.31%  ' }
    @{ Cell = 'D21'; Value = '17.69' }
    @{ Cell = 'E21'; Value = '  +1.00%  ' }
    @{ Cell = 'D22'; Value = '1.003' }
    @{ Cell = 'E22'; Value = '  +0.22%  ' }
    @{ Cell = 'D23'; Value = '12.83' }
    @{ Cell = 'E23'; Value = '  +2.55%  ' }
    @{ Cell = 'D24'; Value = '23.695.70' }
    @{ Cell = 'E24'; Value = '  +1.21%  ' }
    @{ Cell = 'D25'; Value = '2.525' }
    @{ Cell = 'E25'; Value = '  +0.85%  ' }
    @{ Cell = 'D26'; Value = '3.042' }
    @{ Cell = 'E26'; Value = '  -0.62%  ' }
    @{ Cell = 'D27'; Value = '21.28' }
    @{ Cell = 'E27'; Value = '  +0.70%  ' }
    @{ Cell = 'D28'; Value = '151.40' }
    @{ Cell = 'E28'; Value = '  +0.21%  ' }
    @{ Cell = 'E29'; Value = '  -1.00%  ' }
    @{ Cell = 'D30'; Value = '134.39' }
    @{ Cell = 'E30'; Value = '  +1.03%  ' }
    @{ Cell = 'D31'; Value = '1.836.25' }
    @{ Cell = 'E31'; Value = '  +1.32%  ' }
    @{ Cell = 'D32'; Value = '6.901' }
    @{ Cell = 'E32'; Value = '  +4.20%  ' }
    @{ Cell = 'D33'; Value = '2.227' }
    @{ Cell = 'E33'; Value = '  +3.02%  ' }
    @{ Cell = 'D34'; Value = '1.067' }
    @{ Cell = 'E34'; Value = '  +1.69%  ' }
    @{ Cell = 'D35'; Value = '11.67' }
    @{ Cell = 'E35'; Value = '  +2.61%  ' }
    @{ Cell = 'D36'; Value = '0.02804' }
    @{ Cell = 'E36'; Value = '  +1.63%  ' }
    @{ Cell = 'D37'; Value = '0.2510' }
    @{ Cell = 'E37'; Value = '  +0.88%  ' }
    @{ Cell = 'D38'; Value = '0.08827' }
    @{ Cell = 'E38'; Value = '  +0.62%  ' }
    @{ Cell = 'D39'; Value = '6.087' }
    @{ Cell = 'E39'; Value = '  +1.44%  ' }
    @{ Cell = 'D40'; Value = '0.07075' }
    @{ Cell = 'E40'; Value = '  -0.32%  ' }
    @{ Cell = 'D41'; Value = '12.98' }
    @{ Cell = 'E41'; Value = '  +7.29%  ' }
    @{ Cell = 'D42'; Value = '0.7016' }
    @{ Cell = 'E42'; Value = '  +0.30%  ' }
    @{ Cell = 'D43'; Value = '1.339' }
    @{ Cell = 'E43'; Value = '  -0.01%  ' }
    @{ Cell = 'D44'; Value = '16.09' }
    @{ Cell = 'E44'; Value = '  +0.94%  ' }
    @{ Cell = 'D45'; Value = '0.6510' }
    @{ Cell = 'E45'; Value = '  +0.43%  ' }
    @{ Cell = 'D46'; Value = '2.317' }
    @{ Cell = 'E46'; Value = '  +2.06%  ' }
    @{ Cell = 'D47'; Value = '1.000' }
    @{ Cell = 'E47'; Value = '  +0.03%  ' }
    @{ Cell = 'D48'; Value = '3.954' }
    @{ Cell = 'E48'; Value = '  -0.37%  ' }
    @{ Cell = 'D49'; Value = '0.07960' }
    @{ Cell = 'E49'; Value = '  -0.13%  ' }
    @{ Cell = 'D50'; Value = '127.88' }
    @{ Cell = 'E50'; Value = '  +1.10%  ' }
    @{ Cell = 'D51'; Value = '1.185' }
    @{ Cell = 'E51'; Value = '  +0.03%  ' }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $savedStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $u.Value
    $r.Style = $savedStyle
}
